$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers,
# so they remain stored as text (matching original inline-string cells).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "66.766.38"
$ws.Range("E2").Value = "  +5.02%  "
$ws.Range("D3").Value = "3.499.20"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("D5").Value = "592.27"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").Value = "170.47"
$ws.Range("E6").Value = "  +9.25%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.497.62"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  +5.21%  "
$ws.Range("D12").Value = "0.437"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "4.099.96"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "28.07"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").Value = "66.704.38"
$ws.Range("E17").Value = "  +4.76%  "
$ws.Range("D18").Value = "3.499.49"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").Value = "14.07"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "389.12"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("E22").Value = "  +3.24%  "
$ws.Range("D23").Value = "73.01"
$ws.Range("E23").Value = "  +2.88%  "
$ws.Range("D24").Value = "0.998"
$ws.Range("D25").Value = "0.0000125"
$ws.Range("E25").Value = "  +9.43%  "
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").Value = "10.20"
$ws.Range("E27").Value = "  +5.60%  "
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("E31").Value = "  +6.71%  "
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("D33").Value = "23.49"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("D34").Value = "7.40"
$ws.Range("E34").Value = "  +6.55%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("D37").Value = "160.98"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +8.51%  "
$ws.Range("E39").Value = "  +6.18%  "
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("D41").Value = "6.72"
$ws.Range("E41").Value = "  +5.38%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "4.62"
$ws.Range("E42").Value = "  +5.24%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "27.18"
$ws.Range("E43").Value = "  +6.48%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "26.40"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "43.55"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.804.42"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  +12.13%  "
$ws.Range("E48").Value = "  +3.78%  "
$ws.Range("D49").Value = "356.98"
$ws.Range("E49").Value = "  +10.15%  "
$ws.Range("E50").Value = "  +6.65%  "
$ws.Range("D51").Value = "33.00"
$ws.Range("E51").Value = "  +10.58%  "

# Restore default cell style (no explicit style index), matching original file
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
